# Apply betexplorer scrape update (26-10-2023 02:45) to the
# argentina_copa-de-la-liga-profesional_2023 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Swap the match-detail columns (F:V) between row pairs whose
#    fixtures were re-ordered by the scraper (same matchday/date in
#    column E, but the two games swapped position).
# ---------------------------------------------------------------
function Swap-Rows([int]$r1, [int]$r2) {
    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

Swap-Rows 32 33
Swap-Rows 35 36
Swap-Rows 102 103
Swap-Rows 133 134

# ---------------------------------------------------------------
# 2. Append three new fixtures as rows 136-138, copying the
#    formatting of the last existing row (135) so columns A and E
#    keep the correct styles (bold/bordered index, datetime format).
# ---------------------------------------------------------------
function Add-MatchRow(
        [int]$RowNum,
        [int]$Indice,
        [string]$Pais,
        [string]$Torneio,
        [string]$Temporada,
        [double]$DataPartida,
        [string]$Home,
        [int]$HomeGols,
        [string]$Away,
        [int]$AwayGols,
        [double]$HomeOpenOdds,
        [string]$HomeOpenDH,
        [double]$HomeCloseOdds,
        [string]$HomeCloseDH,
        [double]$DrawOpenOdds,
        [string]$DrawOpenDH,
        [double]$DrawCloseOdds,
        [string]$DrawCloseDH,
        [double]$AwayOpenOdds,
        [string]$AwayOpenDH,
        [double]$AwayCloseOdds,
        [string]$AwayCloseDH,
        [string]$Url
    ) {

    $prevRow = $RowNum - 1
    $ws.Range("A$prevRow`:V$prevRow").Copy()
    $ws.Range("A$RowNum`:V$RowNum").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

    $ws.Cells.Item($RowNum, 1).Value = $Indice
    $ws.Cells.Item($RowNum, 2).Value = $Pais
    $ws.Cells.Item($RowNum, 3).Value = $Torneio

    $dCell = $ws.Cells.Item($RowNum, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $Temporada
    $dCell.NumberFormat = "General"

    $ws.Cells.Item($RowNum, 5).Value = $DataPartida
    $ws.Cells.Item($RowNum, 6).Value = $Home
    $ws.Cells.Item($RowNum, 7).Value = $HomeGols
    $ws.Cells.Item($RowNum, 8).Value = $Away
    $ws.Cells.Item($RowNum, 9).Value = $AwayGols
    $ws.Cells.Item($RowNum, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($RowNum, 11).Value = $HomeOpenDH
    $ws.Cells.Item($RowNum, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($RowNum, 13).Value = $HomeCloseDH
    $ws.Cells.Item($RowNum, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($RowNum, 15).Value = $DrawOpenDH
    $ws.Cells.Item($RowNum, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($RowNum, 17).Value = $DrawCloseDH
    $ws.Cells.Item($RowNum, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($RowNum, 19).Value = $AwayOpenDH
    $ws.Cells.Item($RowNum, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($RowNum, 21).Value = $AwayCloseDH
    $ws.Cells.Item($RowNum, 22).Value = $Url
}

Add-MatchRow 136 135 "argentina" "copa-de-la-liga-profesional" "2023" 45224.97916666666 `
    "Atl. Tucuman" 1 "Talleres Cordoba" 0 `
    2.74 "20/10/2023 23:12" 2.55 "25/10/2023 23:27" `
    2.98 "20/10/2023 23:12" 3.06 "25/10/2023 23:27" `
    2.78 "20/10/2023 23:12" 3.18 "25/10/2023 23:27" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/atl-tucuman-talleres-cordoba/jezLm9FN/"

Add-MatchRow 137 136 "argentina" "copa-de-la-liga-profesional" "2023" 45224.97916666666 `
    "Belgrano" 1 "Central Cordoba" 1 `
    1.94 "18/10/2023 22:42" 1.89 "25/10/2023 23:26" `
    3.2 "18/10/2023 22:42" 3.14 "25/10/2023 23:26" `
    4.62 "18/10/2023 22:42" 5.29 "25/10/2023 23:26" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/ca-belgrano-de-cordoba-central-cordoba-santiago-del-estero/Q71j3lFA/"

Add-MatchRow 138 137 "argentina" "copa-de-la-liga-profesional" "2023" 45225.08333333334 `
    "River Plate" 3 "Independiente" 0 `
    1.43 "19/10/2023 22:42" 1.68 "26/10/2023 01:59" `
    4.59 "19/10/2023 22:42" 3.83 "26/10/2023 01:59" `
    7.52 "19/10/2023 22:42" 5.51 "26/10/2023 01:59" `
    "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/river-plate-independiente/AuHirmwo/"

# ---------------------------------------------------------------
# 3. Keep the sheet dimension in sync with the new used range.
# ---------------------------------------------------------------
$ws.UsedRange | Out-Null
